$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 - title slide subtitle: fix "SANDHU" back to "SIDHAU" and re-split
# the runs the way PowerPoint does when the text is retyped.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
$tr1 = $subtitle.TextFrame.TextRange

# Paragraph 1 "CSUMB CST205" (chars 1-12) -> split into "CSUMB " + "CST205"
$tr1.Characters(1, 6).Text = "CSUMB "

# Paragraph 2 starts at char 14 (char 13 is the paragraph mark):
# "PRESENTED BY: DANIEL HOWE, JIWANOT SANDHU, ROGELIO MORENO, NIKOLA PETKOV"
# split into "PRESENTED " + "BY: " + "DANIEL HOWE, JIWANOT SIDHAU, ROGELIO MORENO, NIKOLA PETKOV"
$tr1.Characters(14, 10).Text = "PRESENTED "
$tr1.Characters(24, 4).Text = "BY: "
$tr1.Characters(28, 58).Text = "DANIEL HOWE, JIWANOT SIDHAU, ROGELIO MORENO, NIKOLA PETKOV"

# ---------------------------------------------------------------------------
# Slide 10 - "Image Credits:" textbox: split the trailing colon into its own
# run (matching how PowerPoint splits a run when the colon is retyped).
# ---------------------------------------------------------------------------
$slide10 = $p.Slides.Item(10)
$creditsBox = $slide10.Shapes.Item(2)
$tr10 = $creditsBox.TextFrame.TextRange

# "Image Credits:" occupies chars 23-36 of the full text range; split into
# "Image Credits" (chars 23-35) + ":" (char 36).
$tr10.Characters(23, 13).Text = "Image Credits"
